$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new columns C-I with CIMIS station data
$ws.Range("C1").Value = "x"
$ws.Range("D1").Value = "Downtown"
$ws.Range("E1").Value = "WBAN:23272"
$ws.Range("F1").Value = "Half Moon Bay Airport"
$ws.Range("G1").Value = "WBAN:00228"

$ws.Range("C2").Value = "x"

$ws.Range("C3").Value = "x"
$ws.Range("D3").Value = "SF Bay Reserve"
$ws.Range("E3").Value = "WBAN:99999"

$ws.Range("C4").Value = "x"

$ws.Range("C5").Value = "x"
$ws.Range("D5").Value = "Peraluma Municipal Airport"
$ws.Range("E5").Value = "WBAN:00320"

$ws.Range("C6").Value = "x"
$ws.Range("D6").Value = "Livermore Municipal Airport"
$ws.Range("E6").Value = "WBAN:23285"

$ws.Range("C7").Value = "x"

$ws.Range("C8").Value = "x"
$ws.Range("D8").Value = "SJC"
$ws.Range("E8").Value = "WBAN:23293"
$ws.Range("F8").Value = "Reid-Hillview Airport of Sana Clara"
$ws.Range("G8").Value = "WBAN:93232"

$ws.Range("C9").Value = "x"

$ws.Range("C10").Value = "x"

$ws.Range("C11").Value = "x"
$ws.Range("D11").Value = "SFO"
$ws.Range("E11").Value = "WBAN:23234"
$ws.Range("F11").Value = "San Carlos Airport"
$ws.Range("G11").Value = "WBAN:93231"
$ws.Range("H11").Value = "Palo Alto Airport "
$ws.Range("I11").Value = "WBAN:23289"

$ws.Range("C12").Value = "x"
$ws.Range("D12").Value = "Sonoma County Airport"
$ws.Range("E12").Value = "WBAN:23213"

# Set explicit column widths for the new columns (best-fit sizing)
$ws.Columns.Item(4).ColumnWidth = 25.666666666666668
$ws.Columns.Item(5).ColumnWidth = 11.5
$ws.Columns.Item(6).ColumnWidth = 31.333333333333332
$ws.Columns.Item(7).ColumnWidth = 11.5
$ws.Columns.Item(8).ColumnWidth = 15.5

# Update the active cell selection
$ws.Range("A4").Select() | Out-Null
